$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.594.41'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -2.91%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.632.14'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -3.15%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '590.61'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.90%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '184.00'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.15%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.614'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -3.41%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.17%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.678'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -6.77%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.146'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -10.83%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.77'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -6.22%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -13.78%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.99'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -8.42%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.209.81'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -3.25%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.635.53'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -3.41%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.48%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '67.403.21'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.89%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.46'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -5.83%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.30'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -5.35%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -5.72%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '398.05'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -4.06%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -7.00%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '85.91'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -4.28%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -6.51%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.37'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -4.75%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.52%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.35'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -6.88%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.62'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -9.26%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.04'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -5.64%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.34'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -5.56%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.77'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -8.41%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '67.00'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +2.53%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.93'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -5.11%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '594.03'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -3.47%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -5.69%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '41.84'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -5.71%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.09%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.20%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -7.46%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0₃0739'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -18.85%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -4.00%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.79'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -7.45%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.40'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -13.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.711.42'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -3.03%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -4.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.10'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -3.68%  '
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.55'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -6.78%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'Monero'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '138.22'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -3.30%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.23'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -12.10%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -6.78%  '
